$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in row 10 with the new lab data (Lab 5 - Electron Diffraction)
$ws.Range("A10").Value = "[Wk 10] Monday 21.5.18"
$ws.Range("D10").Value = "Lab Report 4"
$ws.Range("C10").Value = 2
$ws.Range("B10").Value = "1100-1300"

# Update the selection to match the recorded cursor position
$ws.Range("B11").Select()
